$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: new stock item - Trai Oi (Guava)
$ws.Range("A5").Value = "NL006"
$ws.Range("B5").Value = "Gia Vi"
$ws.Range("C5").Value = "Trái Ổi"
$ws.Range("D5").Value = "Ổi sân vườn"
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = "gam"

# Row 6: new stock item - Tom hum alaska (Alaskan lobster)
$ws.Range("A6").Value = "NL007"
$ws.Range("B6").Value = "Hai San"
$ws.Range("C6").Value = "Tôm hùm alaska"
$ws.Range("D6").Value = "lấy ở biển đông"
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = "Kg"
